# Related Work.docx - "added the listing to the theory section"
#
# The _GoBack bookmark (which Word drops at the location of the most
# recent edit) moves from the empty paragraph just above the
# "(Write about your game design)" placeholder to the start of the
# "Here the user is provided a questionnaire ..." paragraph in the
# quiz/theory subsection - i.e. that is where the real edit happened.
# Separately, a space got typed just before the closing ")" of the
# game-design placeholder, splitting its single run into three runs.

$d = $word.ActiveDocument

# --- 1. Move the _GoBack bookmark -----------------------------------

# Remove it from its old spot (empty paragraph right before the
# "(Write about your game design)" placeholder).
$oldMark = $d.Bookmarks("_GoBack")
$oldMark.Delete()

# Re-create it (collapsed) at the start of the paragraph describing the
# quiz questionnaire - the newly-edited "listing" in the theory section.
$quizRange = $d.Content
$quizRange.Find.Execute("Here the user is provided a questionnaire")
$newMarkRange = $d.Range($quizRange.Start, $quizRange.Start)
$d.Bookmarks.Add("_GoBack", $newMarkRange)

# --- 2. Type a space before the closing paren of the placeholder -----

$placeholder = $d.Content
$placeholder.Find.Execute("game design)")
$closeParenStart = $placeholder.End - 1

$caret = $d.Range($closeParenStart, $closeParenStart)
$caret.InsertBefore(" ")

# Nudge formatting on the freshly typed space and back again so it is
# committed as its own run (matching how Word leaves behind separate
# runs around an in-place edit) instead of being re-absorbed into the
# run that precedes it.
$typedSpace = $d.Range($closeParenStart, $closeParenStart + 1)
$typedSpace.Bold = $true

$typedSpaceAgain = $d.Range($closeParenStart, $closeParenStart + 1)
$typedSpaceAgain.Bold = $false

Write-Output "Moved _GoBack bookmark and split the game-design placeholder run."
